$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2323
$ws.Range("I2").Value = 765
$ws.Range("J2").Value = 3361.6667
$ws.Range("K2").Value = 765
$ws.Range("L2").Value = 3361.6667
$ws.Range("M2").Value = -652
$ws.Range("N2").Value = -3587.6667
$ws.Range("H40").Value = 2243.2222
$ws.Range("I40").Value = 2140.1765
$ws.Range("J40").Value = 2418.4
$ws.Range("K40").Value = 2140.1765
$ws.Range("L40").Value = 2418.4
$ws.Range("M40").Value = -1965.1765
$ws.Range("N40").Value = -2768.4
$ws.Range("H74").Value = 98286
$ws.Range("I74").Value = 140986.19
$ws.Range("K74").Value = 140986.19
$ws.Range("M74").Value = -140050.19
$ws.Range("H77").Value = 98286
$ws.Range("I77").Value = 140986.19
$ws.Range("K77").Value = 704930.95
$ws.Range("M77").Value = -700250.95
$ws.Range("H100").Value = 1439.5555
$ws.Range("I100").Value = 650
$ws.Range("J100").Value = 1834.3334
$ws.Range("K100").Value = 650
$ws.Range("L100").Value = 1834.3334
$ws.Range("M100").Value = -109
$ws.Range("N100").Value = -2916.3334
$ws.Range("H125").Value = 150000640
$ws.Range("I125").Value = 142857540
$ws.Range("J125").Value = 166667890
$ws.Range("K125").Value = 1285717860
$ws.Range("L125").Value = 1500011010
$ws.Range("M125").Value = -1285715400
$ws.Range("N125").Value = -1500015930
$ws.Range("H132").Value = 962.43475
$ws.Range("I132").Value = 911.381
$ws.Range("K132").Value = 2734.143
$ws.Range("M132").Value = -204.143

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 800
$ws.Range("I45").Value = 900
$ws.Range("J45").Value = 700
$ws.Range("K45").Value = 900
$ws.Range("L45").Value = 700
$ws.Range("M45").Value = -523
$ws.Range("N45").Value = -1454
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H97").Value = 672.6667
$ws.Range("I97").Value = 701.8570999999999
$ws.Range("J97").Value = 570.5
$ws.Range("K97").Value = 701.8570999999999
$ws.Range("L97").Value = 570.5
$ws.Range("M97").Value = -205.8570999999999
$ws.Range("N97").Value = -1562.5
$ws.Range("H110").Value = 4626342
$ws.Range("J110").Value = 1506.5
$ws.Range("L110").Value = 1506.5
$ws.Range("N110").Value = -5596.5
$ws.Range("H122").Value = 1918.7693
$ws.Range("I122").Value = 1274.5
$ws.Range("J122").Value = 2205.111
$ws.Range("K122").Value = 3823.5
$ws.Range("L122").Value = 6615.333
$ws.Range("M122").Value = -1373.5
$ws.Range("N122").Value = -11515.333
$ws.Range("H134").Value = 74997
$ws.Range("J134").Value = 74997
$ws.Range("L134").Value = 74997
$ws.Range("N134").Value = -85137

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 627.17645
$ws.Range("I22").Value = 460.16666
$ws.Range("K22").Value = 460.16666
$ws.Range("M22").Value = -287.16666
$ws.Range("H86").Value = 773.125
$ws.Range("I86").Value = 791.93335
$ws.Range("J86").Value = 491
$ws.Range("K86").Value = 791.93335
$ws.Range("L86").Value = 491
$ws.Range("M86").Value = 331.06665
$ws.Range("N86").Value = -2737
$ws.Range("H89").Value = 773.125
$ws.Range("I89").Value = 791.93335
$ws.Range("J89").Value = 491
$ws.Range("K89").Value = 3959.66675
$ws.Range("L89").Value = 2455
$ws.Range("M89").Value = 1656.33325
$ws.Range("N89").Value = -13687
$ws.Range("H94").Value = 8537.375
$ws.Range("I94").Value = 9342.714
$ws.Range("J94").Value = 2900
$ws.Range("K94").Value = 9342.714
$ws.Range("L94").Value = 2900
$ws.Range("M94").Value = -8891.714
$ws.Range("N94").Value = -3802

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 28433
$ws.Range("I22").Value = 749.3333
$ws.Range("J22").Value = 42274.832
$ws.Range("K22").Value = 749.3333
$ws.Range("L22").Value = 42274.832
$ws.Range("M22").Value = -399.3333
$ws.Range("N22").Value = -42974.832
$ws.Range("H62").Value = 1999.5
$ws.Range("J62").Value = 1999
$ws.Range("L62").Value = 1999
$ws.Range("N62").Value = -3247
$ws.Range("H65").Value = 1999.5
$ws.Range("J65").Value = 1999
$ws.Range("L65").Value = 9995
$ws.Range("N65").Value = -16235

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H131").Value = 2846.0667
$ws.Range("J131").Value = 3020.7856
$ws.Range("L131").Value = 9062.356800000001
$ws.Range("N131").Value = -19142.3568

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 34399.8
$ws.Range("J15").Value = 35749.75
$ws.Range("L15").Value = 35749.75
$ws.Range("N15").Value = -36325.75
$ws.Range("H43").Value = 8427.666999999999
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H81").Value = 34399.8
$ws.Range("J81").Value = 35749.75
$ws.Range("L81").Value = 35749.75
$ws.Range("N81").Value = -37745.75
$ws.Range("H84").Value = 34399.8
$ws.Range("J84").Value = 35749.75
$ws.Range("L84").Value = 107249.25
$ws.Range("N84").Value = -117233.25
$ws.Range("H97").Value = 251.08333
$ws.Range("I97").Value = 242.09091
$ws.Range("K97").Value = 242.09091
$ws.Range("M97").Value = 253.90909
$ws.Range("H122").Value = 6014.273
$ws.Range("I122").Value = 4807.4287
$ws.Range("K122").Value = 14422.2861
$ws.Range("M122").Value = -11972.2861

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1250
$ws.Range("I22").Value = 1250
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1250
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -955
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1250
$ws.Range("I27").Value = 1250
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1250
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1143
$ws.Range("N27").ClearContents()
$ws.Range("H48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -31322
$ws.Range("H93").Value = 2115.9
$ws.Range("I93").Value = 1837.2858
$ws.Range("J93").Value = 2766
$ws.Range("K93").Value = 1837.2858
$ws.Range("L93").Value = 2766
$ws.Range("M93").Value = -589.2858000000001
$ws.Range("N93").Value = -5262

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 72521.39999999999
$ws.Range("J46").Value = 72521.39999999999
$ws.Range("L46").Value = 72521.39999999999
$ws.Range("N46").Value = -72983.39999999999
$ws.Range("H62").Value = 11062.25
$ws.Range("I62").Value = 7126
$ws.Range("K62").Value = 7126
$ws.Range("M62").Value = -6502
$ws.Range("H65").Value = 11062.25
$ws.Range("I65").Value = 7126
$ws.Range("K65").Value = 35630
$ws.Range("M65").Value = -32510
$ws.Range("H132").Value = 962.3043
$ws.Range("I132").Value = 958.7619
$ws.Range("K132").Value = 2876.2857
$ws.Range("M132").Value = -346.2856999999999
$ws.Range("H134").Value = 72521.39999999999
$ws.Range("J134").Value = 72521.39999999999
$ws.Range("L134").Value = 217564.2
$ws.Range("N134").Value = -222634.2
